$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73; this shifts existing rows 73-86 down to 74-87
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record's data
$ws.Cells.Item(73, 1).Value = 4
$ws.Cells.Item(73, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(73, 3).Value = "Los Lagos"
$ws.Cells.Item(73, 4).Value = 44522
$ws.Cells.Item(73, 5).Value = 10
$ws.Cells.Item(73, 6).Value = 100112022
$ws.Cells.Item(73, 7).Value = "Arveja Verde"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 60
$ws.Cells.Item(73, 11).Value = 20000
$ws.Cells.Item(73, 12).Value = 20000
$ws.Cells.Item(73, 13).Value = 20000
$ws.Cells.Item(73, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(73, 15).Value = "Región del Maule"
$ws.Cells.Item(73, 16).Value = 800
$ws.Cells.Item(73, 17).Value = 25
$ws.Cells.Item(73, 18).Value = "Hortaliza"

# Apply the same date-number format used by the other Fecha cells in column D
$ws.Cells.Item(73, 4).NumberFormat = $ws.Cells.Item(74, 4).NumberFormat
